# Refresh the "Search" sheet's TimeStamp column (C2:C5) with the latest
# test-run values, then let the column width re-fit the (shorter) text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Search")

$ws.Range("C2").Value = "21-09-2020 11:38:0"
$ws.Range("C3").Value = "21-09-2020 11:38:3"
$ws.Range("C4").Value = "21-09-2020 11:38:7"
$ws.Range("C5").Value = "21-09-2020 11:38:10"

# Narrow column C to fit the new (shorter) TimeStamp strings.
$ws.Columns.Item(3).ColumnWidth = 18.8
